$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Gfra1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3912683333333333
$ws.Cells.Item(2, 8).Value = 1.173805
$ws.Cells.Item(2, 9).Value = 0.004923718964983145
$ws.Cells.Item(2, 10).Value = 0.004923718964983145
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.03998533333333333
$ws.Cells.Item(2, 14).Value = 0.119956
$ws.Cells.Item(2, 15).Value = 0.001814551768531471
$ws.Cells.Item(2, 16).Value = 0.00181455176853147
$ws.Cells.Item(2, 17).Value = 0.01564499473111111
$ws.Cells.Item(2, 18).Value = 0.14080495258
$ws.Cells.Item(2, 19).Value = 0.000008934342955662109
$ws.Cells.Item(2, 20).Value = 0.000008934342955662107

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Gfra1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3912683333333333
$ws.Cells.Item(3, 8).Value = 1.173805
$ws.Cells.Item(3, 9).Value = 0.004923718964983145
$ws.Cells.Item(3, 10).Value = 0.004923718964983145
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 15.933008
$ws.Cells.Item(3, 14).Value = 47.799024
$ws.Cells.Item(3, 15).Value = 0.7230468132755195
$ws.Cells.Item(3, 16).Value = 0.7230468132755195
$ws.Cells.Item(3, 17).Value = 6.234081485146667
$ws.Cells.Item(3, 18).Value = 56.10673336632
$ws.Cells.Item(3, 19).Value = 0.003560079307095302
$ws.Cells.Item(3, 20).Value = 0.003560079307095302

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Gfra1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3912683333333333
$ws.Cells.Item(4, 8).Value = 1.173805
$ws.Cells.Item(4, 9).Value = 0.004923718964983145
$ws.Cells.Item(4, 10).Value = 0.004923718964983145
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.062935333333333
$ws.Cells.Item(4, 14).Value = 18.188806
$ws.Cells.Item(4, 15).Value = 0.275138634955949
$ws.Cells.Item(4, 16).Value = 0.275138634955949
$ws.Cells.Item(4, 17).Value = 2.372234602981111
$ws.Cells.Item(4, 18).Value = 21.35011142683
$ws.Cells.Item(4, 19).Value = 0.001354705314932181
$ws.Cells.Item(4, 20).Value = 0.001354705314932181

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Gfra1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.268658666666667
$ws.Cells.Item(5, 8).Value = 18.805976
$ws.Cells.Item(5, 9).Value = 0.07888477275715973
$ws.Cells.Item(5, 10).Value = 0.07888477275715973
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.03998533333333333
$ws.Cells.Item(5, 14).Value = 0.119956
$ws.Cells.Item(5, 15).Value = 0.001814551768531471
$ws.Cells.Item(5, 16).Value = 0.00181455176853147
$ws.Cells.Item(5, 17).Value = 0.2506544063395555
$ws.Cells.Item(5, 18).Value = 2.255889657056
$ws.Cells.Item(5, 19).Value = 0.0001431405039167074
$ws.Cells.Item(5, 20).Value = 0.0001431405039167073

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Gfra1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 6.268658666666667
$ws.Cells.Item(6, 8).Value = 18.805976
$ws.Cells.Item(6, 9).Value = 0.07888477275715973
$ws.Cells.Item(6, 10).Value = 0.07888477275715973
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 15.933008
$ws.Cells.Item(6, 14).Value = 47.799024
$ws.Cells.Item(6, 15).Value = 0.7230468132755195
$ws.Cells.Item(6, 16).Value = 0.7230468132755195
$ws.Cells.Item(6, 17).Value = 99.87858868526935
$ws.Cells.Item(6, 18).Value = 898.9072981674241
$ws.Cells.Item(6, 19).Value = 0.05703738355802786
$ws.Cells.Item(6, 20).Value = 0.05703738355802786

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Gfra1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 6.268658666666667
$ws.Cells.Item(7, 8).Value = 18.805976
$ws.Cells.Item(7, 9).Value = 0.07888477275715973
$ws.Cells.Item(7, 10).Value = 0.07888477275715973
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.062935333333333
$ws.Cells.Item(7, 14).Value = 18.188806
$ws.Cells.Item(7, 15).Value = 0.275138634955949
$ws.Cells.Item(7, 16).Value = 0.275138634955949
$ws.Cells.Item(7, 17).Value = 38.00647212273955
$ws.Cells.Item(7, 18).Value = 342.058249104656
$ws.Cells.Item(7, 19).Value = 0.02170424869521517
$ws.Cells.Item(7, 20).Value = 0.02170424869521517

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Gfra1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 72.68848166666666
$ws.Cells.Item(8, 8).Value = 218.065445
$ws.Cells.Item(8, 9).Value = 0.9147115297293749
$ws.Cells.Item(8, 10).Value = 0.9147115297293749
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.03998533333333333
$ws.Cells.Item(8, 14).Value = 0.119956
$ws.Cells.Item(8, 15).Value = 0.001814551768531471
$ws.Cells.Item(8, 16).Value = 0.00181455176853147
$ws.Cells.Item(8, 17).Value = 2.906473168935555
$ws.Cells.Item(8, 18).Value = 26.15825852042
$ws.Cells.Item(8, 19).Value = 0.001659791423966564
$ws.Cells.Item(8, 20).Value = 0.001659791423966564

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Gfra1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 72.68848166666666
$ws.Cells.Item(9, 8).Value = 218.065445
$ws.Cells.Item(9, 9).Value = 0.9147115297293749
$ws.Cells.Item(9, 10).Value = 0.9147115297293749
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 15.933008
$ws.Cells.Item(9, 14).Value = 47.799024
$ws.Cells.Item(9, 15).Value = 0.7230468132755195
$ws.Cells.Item(9, 16).Value = 0.7230468132755195
$ws.Cells.Item(9, 17).Value = 1158.146159902853
$ws.Cells.Item(9, 18).Value = 10423.31543912568
$ws.Cells.Item(9, 19).Value = 0.6613792566372001
$ws.Cells.Item(9, 20).Value = 0.6613792566372001

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Gfra1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 72.68848166666666
$ws.Cells.Item(10, 8).Value = 218.065445
$ws.Cells.Item(10, 9).Value = 0.9147115297293749
$ws.Cells.Item(10, 10).Value = 0.9147115297293749
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 6.062935333333333
$ws.Cells.Item(10, 14).Value = 18.188806
$ws.Cells.Item(10, 15).Value = 0.275138634955949
$ws.Cells.Item(10, 16).Value = 0.275138634955949
$ws.Cells.Item(10, 17).Value = 440.7055638231855
$ws.Cells.Item(10, 18).Value = 3966.35007440867
$ws.Cells.Item(10, 19).Value = 0.2516724816682082
$ws.Cells.Item(10, 20).Value = 0.2516724816682082

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Ncam1"
$ws.Cells.Item(11, 3).Value = "Gfra1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.117608
$ws.Cells.Item(11, 8).Value = 0.352824
$ws.Cells.Item(11, 9).Value = 0.001479978548482255
$ws.Cells.Item(11, 10).Value = 0.001479978548482255
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.03998533333333333
$ws.Cells.Item(11, 14).Value = 0.119956
$ws.Cells.Item(11, 15).Value = 0.001814551768531471
$ws.Cells.Item(11, 16).Value = 0.00181455176853147
$ws.Cells.Item(11, 17).Value = 0.004702595082666666
$ws.Cells.Item(11, 18).Value = 0.042323355744
$ws.Cells.Item(11, 19).Value = 0.000002685497692537115
$ws.Cells.Item(11, 20).Value = 0.000002685497692537114

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Ncam1"
$ws.Cells.Item(12, 3).Value = "Gfra1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.117608
$ws.Cells.Item(12, 8).Value = 0.352824
$ws.Cells.Item(12, 9).Value = 0.001479978548482255
$ws.Cells.Item(12, 10).Value = 0.001479978548482255
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 15.933008
$ws.Cells.Item(12, 14).Value = 47.799024
$ws.Cells.Item(12, 15).Value = 0.7230468132755195
$ws.Cells.Item(12, 16).Value = 0.7230468132755195
$ws.Cells.Item(12, 17).Value = 1.873849204864
$ws.Cells.Item(12, 18).Value = 16.864642843776
$ws.Cells.Item(12, 19).Value = 0.001070093773196223
$ws.Cells.Item(12, 20).Value = 0.001070093773196223

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Ncam1"
$ws.Cells.Item(13, 3).Value = "Gfra1"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.117608
$ws.Cells.Item(13, 8).Value = 0.352824
$ws.Cells.Item(13, 9).Value = 0.001479978548482255
$ws.Cells.Item(13, 10).Value = 0.001479978548482255
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 6.062935333333333
$ws.Cells.Item(13, 14).Value = 18.188806
$ws.Cells.Item(13, 15).Value = 0.275138634955949
$ws.Cells.Item(13, 16).Value = 0.275138634955949
$ws.Cells.Item(13, 17).Value = 0.7130496986826667
$ws.Cells.Item(13, 18).Value = 6.417447288144
$ws.Cells.Item(13, 19).Value = 0.0004071992775934944
$ws.Cells.Item(13, 20).Value = 0.0004071992775934944

